$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.049.86"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "3.396.46"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'566.09"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.397.71"
$ws.Range("E8").Value = "  +1.10%  "

$ws.Range("E9").Value = "  +1.92%  "

$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("E11").Value = "  +2.23%  "

$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D13").Value = "3.985.63"
$ws.Range("E13").Value = "  +1.32%  "

$ws.Range("E14").Value = "  -3.64%  "

$ws.Range("E15").Value = "  +3.45%  "

$ws.Range("D16").Value = "'26.88"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Value = "63.166.13"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").Value = "3.417.88"
$ws.Range("E18").Value = "  +1.94%  "

$ws.Range("D19").Value = "'6.25"
$ws.Range("E19").Value = "  -4.06%  "

$ws.Range("D20").Value = "'14.02"
$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("D21").Value = "'378.48"
$ws.Range("E21").Value = "  -2.57%  "

$ws.Range("D22").Value = "'8.09"
$ws.Range("E22").Value = "  -4.36%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").Value = "'71.43"
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").Value = "'0.529"
$ws.Range("E25").Value = "  -2.40%  "

$ws.Range("E26").Value = "  +20.36%  "

$ws.Range("D27").Value = "'9.40"
$ws.Range("E27").Value = "  +5.66%  "

$ws.Range("E28").Value = "  -2.48%  "

$ws.Range("E29").Value = "  -0.33%  "

$ws.Range("D30").Value = "'6.01"
$ws.Range("E30").Value = "  +6.11%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("E32").Value = "  +1.52%  "

$ws.Range("D33").Value = "'23.00"
$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'6.36"
$ws.Range("E34").Value = "  -4.34%  "

$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "'6.75"
$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("D37").Value = "'157.62"
$ws.Range("E37").Value = "  -1.81%  "

$ws.Range("E38").Value = "  -1.98%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.921.41"
$ws.Range("E39").Value = "  +3.16%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0758"
$ws.Range("E40").Value = "  +2.18%  "

$ws.Range("E41").Value = "  -4.17%  "

$ws.Range("D42").Value = "'26.78"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'41.08"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.755"
$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("D47").Value = "'23.29"
$ws.Range("E47").Value = "  +4.89%  "

$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("E49").Value = "  +17.56%  "

$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("E51").Value = "  +2.68%  "
